$d = $word.ActiveDocument

# The build regenerated the internal bookmark IDs for these bookmarks
# (pest_table, ind_plots, dv_vs_pred_ipred, prm_vs_iteration) on this
# deploy. Touch each bookmark's ID so the document mints fresh IDs,
# leaving the bookmark names/positions/content untouched.
$names = @("pest_table", "ind_plots", "dv_vs_pred_ipred", "prm_vs_iteration")
foreach ($name in $names) {
    $bm = $d.Bookmarks.Item($name)
    $r = $bm.Range
    $r.BookmarkID = $name
}

Write-Output "done"
